# Updated cryptos list with GitHub Actions.
# Refreshes price (column D) and 1h volume-change (column E) figures, and
# swaps the Litecoin / Cronos rows (28/29) to reflect the new ranking order.
#
# Note: several "Price" values are plain decimal-looking strings
# (e.g. "218.39") that must stay TEXT cells (matching the original
# inlineStr/string cells), not be auto-converted to numbers by Excel's
# "smart" Value assignment. We force that by writing the value with a
# leading apostrophe (Excel's "treat as text" marker) and then resetting
# the cell style to "Normal" so no stray number-format/style sticks to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '87.026.96'
$ws.Range('E2').Value = '  +8.22%  '

$ws.Range('D3').Value = '3.309.13'
$ws.Range('E3').Value = '  +4.78%  '

$ws.Range('E4').Value = '  -0.12%  '

$c = $ws.Range('D5')
$c.Value = '''218.39'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +4.94%  '

$c = $ws.Range('D6')
$c.Value = '''634.33'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.87%  '

$c = $ws.Range('D7')
$c.Value = '''0.323'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +20.08%  '

$ws.Range('E8').Value = '  -0.11%  '

$c = $ws.Range('D9')
$c.Value = '''0.610'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +4.20%  '

$ws.Range('D10').Value = '3.305.71'
$ws.Range('E10').Value = '  +4.14%  '

$ws.Range('E11').Value = '  +2.87%  '

$c = $ws.Range('D12')
$c.Value = '''0.0000271'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +6.50%  '

$ws.Range('E13').Value = '  +2.09%  '

$ws.Range('D14').Value = '3.921.08'
$ws.Range('E14').Value = '  +4.32%  '

$c = $ws.Range('D15')
$c.Value = '''34.39'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +8.84%  '

$c = $ws.Range('D16')
$c.Value = '''5.38'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +2.97%  '

$ws.Range('D17').Value = '86.764.37'
$ws.Range('E17').Value = '  +7.83%  '

$ws.Range('D18').Value = '3.308.73'
$ws.Range('E18').Value = '  +4.06%  '

$c = $ws.Range('D19')
$c.Value = '''14.44'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +2.36%  '

$c = $ws.Range('D20')
$c.Value = '''3.15'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +5.42%  '

$c = $ws.Range('D21')
$c.Value = '''455.31'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +4.63%  '

$c = $ws.Range('D22')
$c.Value = '''9.02'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.85%  '

$c = $ws.Range('D23')
$c.Value = '''5.32'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +3.80%  '

$c = $ws.Range('D24')
$c.Value = '''7.40'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +6.71%  '

$c = $ws.Range('D25')
$c.Value = '''5.35'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +15.43%  '

$c = $ws.Range('D26')
$c.Value = '''12.42'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +15.46%  '

$ws.Range('D27').Value = '3.496.09'
$ws.Range('E27').Value = '  +4.83%  '

$ws.Range('B28').Value = 'Cronos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D28')
$c.Value = '''0.217'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +77.05%  '

$ws.Range('B29').Value = 'Litecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D29')
$c.Value = '''78.10'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +2.84%  '

$ws.Range('E30').Value = '  +6.70%  '

$ws.Range('E31').Value = '  -0.11%  '

$c = $ws.Range('D32')
$c.Value = '''9.22'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +4.31%  '

$c = $ws.Range('D33')
$c.Value = '''590.53'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +6.27%  '

$ws.Range('E34').Value = '  +0.08%  '

$ws.Range('E35').Value = '  +4.92%  '

$c = $ws.Range('D36')
$c.Value = '''2.04'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +3.04%  '

$ws.Range('E37').Value = '  +0.48%  '

$c = $ws.Range('D38')
$c.Value = '''23.36'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +2.52%  '

$c = $ws.Range('D39')
$c.Value = '''6.57'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +18.04%  '

$c = $ws.Range('D40')
$c.Value = '''0.999'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.19%  '

$c = $ws.Range('D41')
$c.Value = '''0.416'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +3.78%  '

$c = $ws.Range('D42')
$c.Value = '''21.41'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.02%  '

$ws.Range('E43').Value = '  +13.87%  '

$c = $ws.Range('D44')
$c.Value = '''3.03'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +13.72%  '

$c = $ws.Range('D45')
$c.Value = '''158.54'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -2.94%  '

$ws.Range('E46').Value = '  +0.05%  '

$c = $ws.Range('D47')
$c.Value = '''188.30'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.36%  '

$c = $ws.Range('D48')
$c.Value = '''46.71'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +9.14%  '

$ws.Range('E49').Value = '  +4.38%  '

$c = $ws.Range('D50')
$c.Value = '''0.782'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.36%  '

$c = $ws.Range('D51')
$c.Value = '''26.29'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +7.39%  '
